$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Occurrences")

# The individualCount column (M) was recorded as the text "?" for every
# occurrence row. Replace it with the actual count (1) as a real number.
for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 13).Value = 1
}

# Append the three UNIPA animal (vertebrate) records that were missing from
# the transcription (rows 34-36).
$rows = @(
    @{ Row = 34; B = "UNCEN-2000NL-HS001-PM001-VE001"; F = "Spilocuscus maculatus" },
    @{ Row = 35; B = "UNCEN-2000NL-HS001-PM001-VE002"; F = "Phalanger orientalis" },
    @{ Row = 36; B = "UNCEN-2000NL-HS001-PM001-VE003"; F = "Phalanger permixtio" }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = "UNCEN-2000NL-HS001-PM001"
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = "Human Observation"
    $ws.Cells.Item($r, 4).Value = "1999-09-11/1999-09-25"
    $ws.Cells.Item($r, 5).Value = "Animalia"
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = "Spesies"
    $ws.Cells.Item($r, 9).Value = "-2.939800"
    $ws.Cells.Item($r, 10).Value = "135.720400"
    $ws.Cells.Item($r, 11).Value = "WGS84"
    $ws.Cells.Item($r, 12).Value = "ID"
    $ws.Cells.Item($r, 13).Value = 1
    $ws.Cells.Item($r, 16).Value = "Present"
}
